$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.628.99"
$ws.Range("E2").Value = "  +0.85%  "

$ws.Range("D3").Value = "2.607.80"
$ws.Range("E3").Value = "  +5.55%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.49%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.604"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.01%  "

$ws.Range("E8").Value = "  +0.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.580"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +13.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.92%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0849"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.97%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +14.35%  "

$ws.Range("D14").Value = "3.002.38"
$ws.Range("E14").Value = "  +5.44%  "

$ws.Range("E15").Value = "  +1.65%  "

$ws.Range("D16").Value = "2.601.10"
$ws.Range("E16").Value = "  +5.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.926"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.71%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.74%  "

$ws.Range("D19").Value = "46.698.62"
$ws.Range("E19").Value = "  +1.03%  "

$ws.Range("E20").Value = "  +7.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "278.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.18%  "

$ws.Range("E25").Value = "  +8.62%  "

$ws.Range("E26").Value = "  +10.94%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "28.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +31.97%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("E29").Value = "  -0.33%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.69%  "

$ws.Range("E31").Value = "  +3.90%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "39.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +15.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.22%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0845"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.37%  "

$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.92%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "152.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.123"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.64%  "

$ws.Range("E40").Value = "  +6.01%  "

$ws.Range("E41").Value = "  +40.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.79%  "

$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.88%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0332"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.33%  "

$ws.Range("D46").Value = "2.117.39"
$ws.Range("E46").Value = "  +6.18%  "

$ws.Range("E47").Value = "  -0.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "93.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.87%  "

$ws.Range("E50").Value = "  +1.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "109.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.90%  "
